# "adding provider to ballon"
# Fill in the new provider_id_number column (X) for every product row in
# the BALON inventory sheet, and tidy up the row heights / selection to
# match the re-saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..10 (row 1 is the header). Each gets the same provider id.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 24).Value = 25998807
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Leave the cursor parked at the bottom-right of the sheet, as in the
# saved workbook, with the view scrolled back to the top-left.
$ws.Range("Z15").Select()
